$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 43 currently has merged cells B43:F43 containing "OK".
# Unmerge and fill C43:F43 with the same value as B43 ("OK").
$ws.Range("B43:F43").UnMerge()
$v = $ws.Range("B43").Value2
$ws.Range("C43").Value = $v
$ws.Range("D43").Value = $v
$ws.Range("E43").Value = $v
$ws.Range("F43").Value = $v

# Update the view state: scroll position and active selection.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("G42").Select()
